$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "NumSeeds" column header to "SeedsWeight" (seed lots now track
# seed weight rather than a raw seed count).
$ws.Range("D1").Value = "SeedsWeight"

# Re-apply the (unchanged) cell style across the used range A1:D10 -- this
# mirrors the formatting touch-up the author made alongside the rename.
$ws.Range("A1:D10").Style = "Normal"

# Move the active selection to D2.
[void]$ws.Range("D2").Select()
